$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "ΜΠΕΝΤΖΑΜΙΝ" / "Ficus benjamina" row (row 58) - it has zero totals
$ws.Rows.Item(58).EntireRow.Delete()

# Remove the "ΦΡΑΓΚΟΣΥΚΙΑ" / "Opuntia ficus-indica" row (originally row 87,
# now row 86 after the first deletion) - it also has zero totals
$ws.Rows.Item(86).EntireRow.Delete()
